# chore: update Sheets via scheduled runner
# Refreshes cached market-price-derived figures (currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ columns H,I,J,K,L,M,N) for specific leve rows
# across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 184.07692
$ws.Range("I11").Value = 184.07692
$ws.Range("K11").Value = 184.07692
$ws.Range("M11").Value = -44.07692

$ws.Range("H17").Value = 5724.577
$ws.Range("I17").Value = 1500
$ws.Range("K17").Value = 4500
$ws.Range("M17").Value = -4332

$ws.Range("H28").Value = 905.6875
$ws.Range("I28").Value = 309.2143
$ws.Range("K28").Value = 309.2143
$ws.Range("M28").Value = 175.7857

$ws.Range("H43").Value = 13137.667
$ws.Range("J43").Value = 13500.5
$ws.Range("L43").Value = 13500.5
$ws.Range("N43").Value = -13638.5

$ws.Range("H51").Value = 7895.2383
$ws.Range("I51").Value = 4700
$ws.Range("J51").Value = 8055
$ws.Range("K51").Value = 4700
$ws.Range("L51").Value = 8055
$ws.Range("M51").Value = -4216
$ws.Range("N51").Value = -9023

$ws.Range("H116").Value = 8221.143
$ws.Range("I116").Value = 6473.25
$ws.Range("J116").Value = 10551.667
$ws.Range("K116").Value = 6473.25
$ws.Range("L116").Value = 10551.667
$ws.Range("M116").Value = -3031.25
$ws.Range("N116").Value = -17435.667

$ws.Range("H132").Value = 704.1177
$ws.Range("I132").Value = 623.125
$ws.Range("K132").Value = 1869.375
$ws.Range("M132").Value = 660.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2798.611
$ws.Range("I2").Value = 1323.2142
$ws.Range("K2").Value = 1323.2142
$ws.Range("M2").Value = -1210.2142

$ws.Range("H97").Value = 1823.8
$ws.Range("I97").Value = 1823.8
$ws.Range("K97").Value = 1823.8
$ws.Range("M97").Value = -1327.8

$ws.Range("H102").Value = 2593.6428
$ws.Range("I102").Value = 1181.7
$ws.Range("K102").Value = 1181.7
$ws.Range("M102").Value = 440.3

$ws.Range("H116").Value = 2798.611
$ws.Range("I116").Value = 1323.2142
$ws.Range("K116").Value = 1323.2142
$ws.Range("M116").Value = 970.7858000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2798.611
$ws.Range("I3").Value = 1323.2142
$ws.Range("K3").Value = 1323.2142
$ws.Range("M3").Value = -1209.2142

$ws.Range("H99").Value = 2609.75
$ws.Range("I99").Value = 2117.0667
$ws.Range("K99").Value = 2117.0667
$ws.Range("M99").Value = -619.0666999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3892.7917
$ws.Range("I58").Value = 1949.9231
$ws.Range("K58").Value = 1949.9231
$ws.Range("M58").Value = -1746.9231

$ws.Range("H122").Value = 6020.1333
$ws.Range("I122").Value = 2163.9092
$ws.Range("K122").Value = 6491.7276
$ws.Range("M122").Value = -4041.7276

$ws.Range("H132").Value = 3359.12
$ws.Range("I132").Value = 2398.2
$ws.Range("K132").Value = 7194.599999999999
$ws.Range("M132").Value = -4664.599999999999

$ws.Range("H133").Value = 36882.93
$ws.Range("J133").Value = 38196.75
$ws.Range("L133").Value = 38196.75
$ws.Range("N133").Value = -43256.75

$ws.Range("H134").Value = 4113.6875
$ws.Range("I134").Value = 2977.8
$ws.Range("J134").Value = 6006.8335
$ws.Range("K134").Value = 8933.400000000001
$ws.Range("L134").Value = 18020.5005
$ws.Range("M134").Value = -6398.400000000001
$ws.Range("N134").Value = -23090.5005

$ws.Range("H136").Value = 3892.7917
$ws.Range("I136").Value = 1949.9231
$ws.Range("K136").Value = 5849.7693
$ws.Range("M136").Value = -3299.7693

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 287858.16
$ws.Range("J37").Value = 287858.16
$ws.Range("L37").Value = 863574.48
$ws.Range("N37").Value = -863798.48

$ws.Range("H132").Value = 3442.926
$ws.Range("J132").Value = 4928.8
$ws.Range("L132").Value = 44359.2
$ws.Range("N132").Value = -49419.2

$ws.Range("H137").Value = 65355.5
$ws.Range("J137").Value = 74527.71000000001
$ws.Range("L137").Value = 223583.13
$ws.Range("N137").Value = -233783.13

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1532.8889
$ws.Range("I2").Value = 297.66666
$ws.Range("J2").Value = 4003.3333
$ws.Range("K2").Value = 297.66666
$ws.Range("L2").Value = 4003.3333
$ws.Range("M2").Value = -184.66666
$ws.Range("N2").Value = -4229.3333

$ws.Range("H123").Value = 37628.2
$ws.Range("J123").Value = 37628.2
$ws.Range("L123").Value = 37628.2
$ws.Range("N123").Value = -42528.2

$ws.Range("H126").Value = 3088.862
$ws.Range("I126").Value = 1610.4667
$ws.Range("K126").Value = 4831.4001
$ws.Range("M126").Value = -2361.4001

$ws.Range("H132").Value = 3173.4644
$ws.Range("I132").Value = 2602.8696
$ws.Range("K132").Value = 7808.6088
$ws.Range("M132").Value = -5278.6088

$ws.Range("H134").Value = 77777
$ws.Range("J134").Value = 77777
$ws.Range("L134").Value = 233331
$ws.Range("N134").Value = -238401

$ws.Range("H135").Value = 69223
$ws.Range("J135").Value = 69223
$ws.Range("L135").Value = 69223
$ws.Range("N135").Value = -79363

$ws.Range("H136").Value = 21521.5
$ws.Range("J136").Value = 19654.21
$ws.Range("L136").Value = 58962.63
$ws.Range("N136").Value = -64062.63

$ws.Range("H140").Value = 59930
$ws.Range("J140").Value = 59930
$ws.Range("L140").Value = 59930
$ws.Range("N140").Value = -70290

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3876.6
$ws.Range("I22").Value = 1034.8
$ws.Range("J22").Value = 5297.5
$ws.Range("K22").Value = 1034.8
$ws.Range("L22").Value = 5297.5
$ws.Range("M22").Value = -739.8
$ws.Range("N22").Value = -5887.5

$ws.Range("H27").Value = 3876.6
$ws.Range("I27").Value = 1034.8
$ws.Range("J27").Value = 5297.5
$ws.Range("K27").Value = 1034.8
$ws.Range("L27").Value = 5297.5
$ws.Range("M27").Value = -927.8
$ws.Range("N27").Value = -5511.5

$ws.Range("H55").Value = 1788226.5
$ws.Range("I55").Value = 3334772.8
$ws.Range("J55").Value = 3750.077
$ws.Range("K55").Value = 3334772.8
$ws.Range("L55").Value = 3750.077
$ws.Range("M55").Value = -3334599.8
$ws.Range("N55").Value = -4096.077

$ws.Range("H93").Value = 2264.1738
$ws.Range("I93").Value = 2272.762
$ws.Range("K93").Value = 2272.762
$ws.Range("M93").Value = -1024.762

$ws.Range("H100").Value = 18102
$ws.Range("I100").Value = 18000
$ws.Range("J100").Value = 18170
$ws.Range("K100").Value = 18000
$ws.Range("L100").Value = 18170
$ws.Range("M100").Value = -17459
$ws.Range("N100").Value = -19252

$ws.Range("H132").Value = 3207.0256
$ws.Range("I132").Value = 3015.52
$ws.Range("J132").Value = 3549
$ws.Range("K132").Value = 9046.559999999999
$ws.Range("L132").Value = 10647
$ws.Range("M132").Value = -6516.559999999999
$ws.Range("N132").Value = -15707

$ws.Range("H141").Value = 69999
$ws.Range("J141").Value = 69999
$ws.Range("L141").Value = 69999
$ws.Range("N141").Value = -80359

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3422.4443
$ws.Range("I122").Value = 799.5
$ws.Range("J122").Value = 8668.333000000001
$ws.Range("K122").Value = 2398.5
$ws.Range("L122").Value = 26004.999
$ws.Range("M122").Value = 51.5
$ws.Range("N122").Value = -30904.999

$ws.Range("H130").Value = 71969
$ws.Range("J130").Value = 71969
$ws.Range("L130").Value = 71969
$ws.Range("N130").Value = -82009
